$wb = $excel.ActiveWorkbook

$msg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b41a868f3a8d4acdc6708b8048446e1265b69681/e2e/c164da8b-4831-46da-82c3-894f42412a10.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/60c9c1e05b3499035788465f058567a6740f1d0e/e2e/c164da8b-4831-46da-82c3-894f42412a10.md."

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Columns.Item(16).ColumnWidth = 39.17

$targetMd = "c164da8b-4831-46da-82c3-894f42412a10.md"
$targetXlfZh = "c164da8b-4831-46da-82c3-894f42412a10.6890eb99f564a45e6c02c31b531789f0a21b50c9.zh-cn.xlf"
$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/60c9c1e05b3499035788465f058567a6740f1d0e/e2e/c164da8b-4831-46da-82c3-894f42412a10.md"

$cellI7 = $wsZh.Range("I7")
$cellI7.Style = "Hyperlink"
$cellI7.Font.Underline = 2
$cellI7.Font.Color = 15570276
$cellI7.Value = $targetMd
$wsZh.Hyperlinks.Add($cellI7, $latestUrl, [Type]::Missing, [Type]::Missing, $targetMd) | Out-Null

$wsZh.Range("J7").Value = $targetXlfZh
$wsZh.Range("K7").Value = "2016-08-25 10:44:57"
$wsZh.Range("P7").Value = $msg

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(16).ColumnWidth = 39.17

$targetXlfDe = "c164da8b-4831-46da-82c3-894f42412a10.6890eb99f564a45e6c02c31b531789f0a21b50c9.de-de.xlf"

$cellI7de = $wsDe.Range("I7")
$cellI7de.Style = "Hyperlink"
$cellI7de.Font.Underline = 2
$cellI7de.Font.Color = 15570276
$cellI7de.Value = $targetMd
$wsDe.Hyperlinks.Add($cellI7de, $latestUrl, [Type]::Missing, [Type]::Missing, $targetMd) | Out-Null

$wsDe.Range("J7").Value = $targetXlfDe
$wsDe.Range("K7").Value = "2016-08-25 10:45:17"
$wsDe.Range("P7").Value = $msg

Write-Host "done"
